# "Start depth image invisible" - reposition/resize the depth-map picture
# on the "Look At" (slide 4) and "Don't Look" (slide 5) slides so that it
# no longer lines up with the color photo at the start (i.e. starts the
# reveal "invisible"/mismatched).
#
# Note: Shape.Left/Top/Width/Height round-trip through a 32-bit "points"
# value in this COM host, so the literals below are chosen to be the
# exact point values whose float32 representation converts back to the
# target EMU (value/12700) instead of the naive double division, which
# can drift the stored EMU by 1.

$p = $ppt.ActivePresentation

# --- Slide 4 ("Far & Look At") -> Picture 3 (the depth/left image) ---
$s4 = $p.Slides.Item(4)
$pic4 = $s4.Shapes.Item("Picture 3")
$pic4.Left   = 30.0
$pic4.Top    = 204.0
$pic4.Width  = 362.5644226074219
$pic4.Height = 171.98574829101562

# --- Slide 5 ("Far & Don't Look") -> Picture 1 (the depth/left image) ---
$s5 = $p.Slides.Item(5)
$pic5 = $s5.Shapes.Item("Picture 1")
$pic5.Left   = 48.0
$pic5.Top    = 210.0
$pic5.Width  = 362.6276550292969
$pic5.Height = 171.98574829101562
